$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unmerge B5:F5 first
$ws.Range("B5:F5").UnMerge()

# Copy fill formatting from E4 (fillId4/theme5) onto E5, and F4 (fillId5/theme9) onto F5
$ws.Range("E4").Copy() | Out-Null
$ws.Range("E5").PasteSpecial(-4122) | Out-Null
$ws.Range("F4").Copy() | Out-Null
$ws.Range("F5").PasteSpecial(-4122) | Out-Null

# Set F5 new text value
$ws.Range("F5").Value2 = '$document'

# Re-merge B5:D5
$ws.Range("B5:D5").Merge()

# New cell E6 takes fillId4 like row6 B6:D6
$ws.Range("D6").Copy() | Out-Null
$ws.Range("E6").PasteSpecial(-4122) | Out-Null
# F6 should carry fillId5 like F4
$ws.Range("F4").Copy() | Out-Null
$ws.Range("F6").PasteSpecial(-4122) | Out-Null
# New cell E7 takes fillId7 like row7 B7:D7
$ws.Range("D7").Copy() | Out-Null
$ws.Range("E7").PasteSpecial(-4122) | Out-Null

# Apply thin border to table body A3:F9
$ws.Range("A3:F9").Borders.LineStyle = 1

Write-Output "done"
